$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent results for the 380 kV case (rows 2-25, columns C-H, L, N)
$colC = @(
    13.55900185796705,
    13.4816248690508,
    13.43837151218974,
    13.4218211190916,
    13.41913803068623,
    13.43814394651769,
    13.5314426164512,
    13.74785763166715,
    13.9267231413649,
    14.01226241521013,
    14.04523948718304,
    14.03811150536243,
    14.01496381921019,
    14.00086094517313,
    13.92121567229163,
    13.87341440701757,
    13.84631408193672,
    13.83720642639589,
    13.87846229037131,
    14.02174709709875,
    14.11879537796627,
    14.06669265417337,
    13.87617895366222,
    13.68577831548532
)
for ($i = 0; $i -lt $colC.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}

$colD = @(
    8.668425739847907,
    8.600279589338216,
    8.561283301186396,
    8.54611529902918,
    8.543640545594982,
    8.561075802476182,
    8.644343395871864,
    8.829803992956032,
    8.978948137811136,
    9.049437807732847,
    9.076495748221735,
    9.070652349051427,
    9.051656653790264,
    9.0400683206179,
    8.974393390303426,
    8.934770471476583,
    8.912230166716759,
    8.904641748294285,
    8.938962657424467,
    9.057226366662372,
    9.136638772420456,
    9.094066027348282,
    8.937066623223368,
    8.777318547159412
)
for ($i = 0; $i -lt $colD.Count; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}

$colE = @(
    13.94350578747531,
    13.90519190635668,
    13.88584101810917,
    13.87900244070478,
    13.87793004497138,
    13.88574455525301,
    13.92942766219119,
    14.04827398022687,
    14.15584331716709,
    14.20915972431979,
    14.22997613552102,
    14.22546514605383,
    14.21085976804016,
    14.20199503713991,
    14.15244687154543,
    14.12317094580974,
    14.1067449960712,
    14.10125444054148,
    14.12624471867985,
    14.21513274800013,
    14.27687643601536,
    14.24359010420323,
    14.12485380369239,
    14.01256355229876
)
for ($i = 0; $i -lt $colE.Count; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $colE[$i]
}

$colF = @(
    47.9076808159054,
    46.83405872097018,
    46.18252198220152,
    45.91932725982485,
    45.87577481330221,
    46.17896257887734,
    47.53613242312501,
    50.24183648997253,
    52.2352237173999,
    53.13899742026904,
    53.48046691368538,
    53.40696469294019,
    53.16710759304178,
    53.02007835300726,
    52.17606996575329,
    51.65727789742891,
    51.3586161282243,
    51.25745918716812,
    51.71253424990176,
    53.23758289492702,
    54.22966244685191,
    53.7006993027741,
    51.68755407727866,
    49.50739887925336
)
for ($i = 0; $i -lt $colF.Count; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $colF[$i]
}

$colG = @(
    72.18944697254682,
    70.0903229290731,
    68.80107669441772,
    68.27646305094002,
    68.18942332164092,
    68.79399731418533,
    71.46625081769692,
    76.66788407378181,
    80.42167366441181,
    82.10636623692731,
    82.74040831249853,
    82.60403931071448,
    82.15861103951644,
    81.88524626135931,
    80.31105956290855,
    79.33901692144808,
    78.77779905419064,
    78.58743521373179,
    79.44271707698709,
    82.28955498714295,
    84.12706683342137,
    83.14865113589286,
    79.39584159484056,
    75.26962022524248
)
for ($i = 0; $i -lt $colG.Count; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $colG[$i]
}

$colH = @(
    23.23254646083079,
    22.94706004378695,
    22.77866971492673,
    22.71184656557442,
    22.70086091555697,
    22.77776115841678,
    23.13271749797991,
    23.88059320178671,
    24.45749077920002,
    24.72496789563602,
    24.82690076231405,
    24.80492027918215,
    24.73334150394997,
    24.68957907039972,
    24.4401051549455,
    24.28829572514255,
    24.20145849695763,
    24.17214170793591,
    24.30440704848571,
    24.75434905562831,
    25.05213535015073,
    24.89288692482633,
    24.29712174250673,
    23.67314524421316
)
for ($i = 0; $i -lt $colH.Count; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $colH[$i]
}

$colL = @(
    9.315572113616621,
    9.332952975339234,
    9.345659042175612,
    9.351345800527138,
    9.352320744637565,
    9.345733678553058,
    9.321141313725869,
    9.289170125078467,
    9.275745737809839,
    9.271856640215633,
    9.27070539387131,
    9.270938997608424,
    9.271755469235908,
    9.272297521850428,
    9.276044746176822,
    9.278913400493986,
    9.280771967287276,
    9.281436993052242,
    9.278586420053449,
    9.271506906931011,
    9.268754965118202,
    9.270051338267285,
    9.278733596060148,
    9.296062750898752
)
for ($i = 0; $i -lt $colL.Count; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $colL[$i]
}

$colN = @(
    18.99769123883691,
    18.40031508502705,
    18.02485520896364,
    17.86990355188116,
    17.84406337566596,
    18.02277304766463,
    18.79364780656866,
    20.2273683202997,
    21.22223697909767,
    21.66018057919901,
    21.82377585681918,
    21.78864458690801,
    21.67368539489661,
    21.60297336126124,
    21.19330956972086,
    20.93814219015166,
    20.79000725568142,
    20.73962067985786,
    20.96544799483449,
    21.70751365554063,
    22.17935961385674,
    21.92877110911182,
    20.95310750188673,
    19.84905939529497
)
for ($i = 0; $i -lt $colN.Count; $i++) {
    $ws.Cells.Item($i + 2, 14).Value = $colN[$i]
}
